# Add 2022-Q4 data
# 1. "总计" (Total) sheet: insert a new row right under the header with the
#    2022-Q4 summary (date / holding count / holding market value) and shift
#    the existing quarters down by one row, renumbering the index column.
# 2. Insert a brand-new worksheet named "2022-Q4" right after "总计" holding
#    the per-fund detail rows for that quarter (mirrors the layout already
#    used by the other quarter sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet - insert the 2022-Q4 summary row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 2.08

# Match the look of the other index-column cells (bold / bordered / centered)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Renumber the index column for the rows that shifted down one slot
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6

# ---------------------------------------------------------------------
# 2. New "2022-Q4" sheet - per-fund holding detail
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$header = $q4.Range("B1:H1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'519692"
$q4.Range("C2").Value = "交银成长混合A"
$q4.Range("D2").Value = "'25.08"
$q4.Range("E2").Value = "'80.90"
$q4.Range("F2").Value = "'3.44"
$q4.Range("G2").Value = "'0.8628"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'519694"
$q4.Range("C3").Value = "交银蓝筹混合"
$q4.Range("D3").Value = "'17.19"
$q4.Range("E3").Value = "'81.14"
$q4.Range("F3").Value = "'3.30"
$q4.Range("G3").Value = "'0.5673"
$q4.Range("H3").Value = 9

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'011184"
$q4.Range("C4").Value = "东方阿尔法招阳混合A"
$q4.Range("D4").Value = "'6.19"
$q4.Range("E4").Value = "'94.10"
$q4.Range("F4").Value = "'9.14"
$q4.Range("G4").Value = "'0.5658"
$q4.Range("H4").Value = 3

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'011185"
$q4.Range("C5").Value = "东方阿尔法招阳混合C"
$q4.Range("D5").Value = "'0.88"
$q4.Range("E5").Value = "'94.10"
$q4.Range("F5").Value = "'9.14"
$q4.Range("G5").Value = "'0.0804"
$q4.Range("H5").Value = 3

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'960016"
$q4.Range("C6").Value = "交银成长混合H"
$q4.Range("D6").Value = "'0.17"
$q4.Range("E6").Value = "'80.90"
$q4.Range("F6").Value = "'3.44"
$q4.Range("G6").Value = "'0.0058"
$q4.Range("H6").Value = 8

$idxCol = $q4.Range("A2:A6")
$idxCol.Font.Bold = $true
$idxCol.Borders.LineStyle = 1
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160

Write-Output "2022-Q4 sheet added and 总计 updated"
